$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the visitor record in row 2 (A=email, B=firstName, C=lastName, D=phone)
# replacing Sam Wesh's details with Paul Murimi's details.
$ws.Range("B2").Value = "Paul"
$ws.Range("C2").Value = "Murimi"
$ws.Range("A2").Value = "paulmaina@gmail.com"
$ws.Range("D2").Value = "0725165221"
